$d = $word.ActiveDocument

# --- French typographic fix: regular space -> non-breaking space before ":" ---
$d.Content.Find.Execute("Organisation : Fabrikam Inc.", $true, $false, $false, $false, $false, $true, 1, $false, "Organisation : Fabrikam Inc.", 2)
$d.Content.Find.Execute("Date : 27 juin 2024", $true, $false, $false, $false, $false, $true, 1, $false, "Date : 27 juin 2024", 2)

# --- Wording / copy-edit changes ---
$d.Content.Find.Execute("Fabrikam Inc. a fait l’objet d’une initiative complète de transformation numérique visant à améliorer l’efficacité opérationnelle, à améliorer l’expérience client et à stimuler l’innovation.", $true, $false, $false, $false, $false, $true, 1, $false, "Fabrikam Inc. a mené une initiative de transformation numérique complète visant à améliorer son efficacité opérationnelle, à améliorer l’expérience client et à stimuler l’innovation.", 2)
$d.Content.Find.Execute("Le résumé suivant décrit les principales mises à jour et les jalons réalisés à ce jour.", $true, $false, $false, $false, $false, $true, 1, $false, "Le résumé suivant décrit les principales mises à jour et les étapes majeures réalisées à ce jour.", 2)
$d.Content.Find.Execute("Mises à jour clés", $true, $false, $false, $false, $false, $true, 1, $false, "Mises à jour principales", 2)
$d.Content.Find.Execute("Implémentation de l’infrastructure cloud", $true, $false, $false, $false, $false, $true, 1, $false, "Implémentation d’une infrastructure cloud", 2)
$d.Content.Find.Execute("Migration de 80 % des applications locales vers le cloud.", $true, $false, $false, $false, $false, $true, 1, $false, "Migration de 80 % des applications locales vers le cloud", 2)
$d.Content.Find.Execute("Amélioration de la scalabilité et réduction des coûts informatiques de 25 %.", $true, $false, $false, $false, $false, $true, 1, $false, "Amélioration de la scalabilité et réduction des coûts informatiques de 25 %", 2)
$d.Content.Find.Execute("Sécurité et conformité des données améliorées avec les normes du secteur.", $true, $false, $false, $false, $false, $true, 1, $false, "Amélioration de la sécurité et de la conformité des données avec les normes du secteur", 2)
$d.Content.Find.Execute("Intégration d’analyses basées sur l’IA pour simplifier les processus décisionnels.", $true, $false, $false, $false, $false, $true, 1, $false, "Intégration d’analyses basées sur l’IA pour simplifier les processus décisionnels", 2)
$d.Content.Find.Execute("Déploiement de modèles Machine Learning pour prédire le comportement des clients et personnaliser les efforts marketing.", $true, $false, $false, $false, $false, $true, 1, $false, "Déploiement de modèles Machine Learning pour prédire le comportement des clients et personnaliser les initiatives marketing", 2)
$d.Content.Find.Execute("Réduction des processus manuels, ce qui entraîne une augmentation de 30 % de la productivité.", $true, $false, $false, $false, $false, $true, 1, $false, "Réduction des processus manuels, entraînant une augmentation de 30 % de la productivité", 2)
$d.Content.Find.Execute("Lancement d’un nouveau portail client avec des fonctionnalités en libre-service.", $true, $false, $false, $false, $false, $true, 1, $false, "Lancement d’un nouveau portail client avec des fonctionnalités en libre-service", 2)
$d.Content.Find.Execute("Introduction des chatbots pour le support client 24/7, ce qui réduit les temps de réponse de 50 %.", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction de chatbots pour un support client 24/7, réduisant les temps de réponse de 50 %", 2)
$d.Content.Find.Execute("Amélioration de la satisfaction des clients de 20 % au cours de l’année dernière.", $true, $false, $false, $false, $false, $true, 1, $false, "Amélioration de la satisfaction des clients de 20 % au cours de l’année dernière", 2)
$d.Content.Find.Execute("Implémentation de l’automatisation des processus robotisés (RPA) pour les tâches courantes.", $true, $false, $false, $false, $false, $true, 1, $false, "Implémentation de l’automatisation robotisée des processus (RPA) pour les tâches courantes", 2)
$d.Content.Find.Execute("A atteint une réduction de 40 % du temps de traitement pour les opérations commerciales clés.", $true, $false, $false, $false, $false, $true, 1, $false, "Réduction de 40 % du temps de traitement pour les opérations clés de l’entreprise", 2)
$d.Content.Find.Execute("Réaffectation des ressources humaines à des rôles plus stratégiques au sein de l’organisation.", $true, $false, $false, $false, $false, $true, 1, $false, "Réaffectation des ressources humaines à des rôles plus stratégiques au sein de l’organisation", 2)
$d.Content.Find.Execute("Des programmes d’alphabétisation numérique ont été menés pour tous les employés.", $true, $false, $false, $false, $false, $true, 1, $false, "Programmes de formation numérique pour tous les employés", 2)
$d.Content.Find.Execute("A lancé une nouvelle plateforme d’apprentissage électronique avec des cours sur les technologies émergentes.", $true, $false, $false, $false, $false, $true, 1, $false, "Nouvelle plateforme de formation en ligne avec des cours sur les technologies émergentes", 2)
$d.Content.Find.Execute("Augmentation de l’engagement des employés et de l’adoption de nouveaux outils de 35 %.", $true, $false, $false, $false, $false, $true, 1, $false, "Augmentation de 35 %. de l’engagement des employés et de l’adoption de nouveaux outils", 2)
$d.Content.Find.Execute("Q1 2024 : Migration terminée vers l’infrastructure cloud.", $true, $false, $false, $false, $false, $true, 1, $false, "T1 2024 : Migration vers l’infrastructure cloud", 2)
$d.Content.Find.Execute("Q2 2024 : Plateforme d’analytique basée sur l’IA lancée.", $true, $false, $false, $false, $false, $true, 1, $false, "T2 2024 : Lancement de la plateforme d’analytique basée sur l’IA", 2)
$d.Content.Find.Execute("Q3 2024 : Introduction du nouveau portail client numérique.", $true, $false, $false, $false, $false, $true, 1, $false, "T3 2024 : Introduction du nouveau portail client numérique", 2)
$d.Content.Find.Execute("Q4 2024 : 50 % d’automatisation des processus de routine.", $true, $false, $false, $false, $false, $true, 1, $false, "T4 2024 : 50 % d’automatisation des processus courants", 2)
$d.Content.Find.Execute("Poursuivez l’expansion des applications IA et Machine Learning dans tous les services.", $true, $false, $false, $false, $false, $true, 1, $false, "Poursuivre l’expansion des applications IA et Machine Learning dans tous les services", 2)
$d.Content.Find.Execute("Améliorez davantage l’expérience client numérique avec de nouvelles fonctionnalités et services.", $true, $false, $false, $false, $false, $true, 1, $false, "Améliorer davantage l’expérience client numérique avec de nouvelles fonctionnalités et de nouveaux services", 2)
$d.Content.Find.Execute("Concentrez-vous sur les mesures de cybersécurité pour vous protéger contre les menaces en constante évolution.", $true, $false, $false, $false, $false, $true, 1, $false, "Se concentrer sur les mesures de cybersécurité pour se protéger contre les menaces en constante évolution", 2)
$d.Content.Find.Execute("Développez une stratégie numérique complète pour les cinq prochaines années.", $true, $false, $false, $false, $false, $true, 1, $false, "Développer une stratégie numérique complète pour les cinq prochaines années", 2)
$d.Content.Find.Execute("L’organisation reste engagée à tirer parti de la technologie pour stimuler la croissance et l’innovation futures.", $true, $false, $false, $false, $false, $true, 1, $false, "L’organisation reste déterminée à tirer parti de la technologie pour stimuler la croissance et l’innovation.", 2)

# --- Bold the two headings that were inconsistently left un-bolded ---
$r = $d.Content
$r.Find.Execute("Automatisation des processus", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = 1

$r2 = $d.Content
$r2.Find.Execute("Plans pour l’avenir", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Font.Bold = 1
